# Auto-generated script to add sheet "2025-12-10" with weekly ranking data
$wb = $excel.ActiveWorkbook
$srcSheet = $wb.Worksheets.Item("2025-12-03")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "2025-12-10"

# Header row: copy format (style index 1: bold, border, centered) from existing sheet, then set values
$srcSheet.Range("A1:D1").Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)
$ws.Range("A1").Value = "rank"
$ws.Range("B1").Value = "title"
$ws.Range("C1").Value = "volume"
$ws.Range("D1").Value = "publisher"

# Data rows
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "薬屋のひとりごと"
$ws.Cells.Item(2, 3).Value = 16
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = "ブルーロック"
$ws.Cells.Item(3, 3).Value = 36
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "ぼっち・ざ・ろっく!"
$ws.Cells.Item(4, 3).Value = 8
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "ハズレ枠ので最強になった俺がすべてを蹂躙するまで"
$ws.Cells.Item(5, 3).Value = 13
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = "ブルーピリオド"
$ws.Cells.Item(6, 3).Value = 18
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "黒執事"
$ws.Cells.Item(7, 3).Value = 35
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "異世界でスローライフを(願望)"
$ws.Cells.Item(8, 3).Value = 1
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "GIANT KILLING"
$ws.Cells.Item(9, 3).Value = 68
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "ワンパンマン"
$ws.Cells.Item(10, 3).Value = 35
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "ハズレ枠ので最強になった俺がすべてを蹂躙するまで"
$ws.Cells.Item(11, 3).Value = 1
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "Sランク冒険者である俺の娘たちは重度のファザコンでした"
$ws.Cells.Item(12, 3).Value = 1
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = "神血の救世主~0.00000001%を引き当て最強へ~"
$ws.Cells.Item(13, 3).Value = 6
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = "勇者パーティーを追放された白魔導師、Sランク冒険者に拾われる~この白魔導師が規格外すぎる~(コミック)"
$ws.Cells.Item(14, 3).Value = 10
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = "MIX"
$ws.Cells.Item(15, 3).Value = 24
$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = "オシバナ!"
$ws.Cells.Item(16, 3).Value = 1
$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = "元婚約者から逃げるため吸血伯爵に恋人のフリをお願いしたら、なぜか溺愛モードになりました"
$ws.Cells.Item(17, 3).Value = 6
$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(18, 2).Value = "SPY×FAMILY"
$ws.Cells.Item(18, 3).Value = 16
$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(19, 2).Value = "社内探偵"
$ws.Cells.Item(19, 3).Value = 66
$ws.Cells.Item(20, 1).Value = 19
$ws.Cells.Item(20, 2).Value = "BLOOD THE LAST VAMPIRE 2000"
$ws.Cells.Item(20, 3).Value = 1
$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(21, 2).Value = "クトゥルーとか全然わからない俺が、邪神の力で爆乳女子と無双する1"
$ws.Cells.Item(21, 3).Value = 1
$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(22, 2).Value = "1・2のアッホ!!"
$ws.Cells.Item(22, 3).Value = 1
$ws.Cells.Item(23, 1).Value = 22
$ws.Cells.Item(23, 2).Value = "×××HOLiC・戻"
$ws.Cells.Item(23, 3).Value = 5
$ws.Cells.Item(24, 1).Value = 23
$ws.Cells.Item(24, 2).Value = "死に戻りの幸薄令嬢、今世では最恐ラスボスお義兄様に溺愛されてます"
$ws.Cells.Item(24, 3).Value = 12
$ws.Cells.Item(25, 1).Value = 24
$ws.Cells.Item(25, 2).Value = "勇者殺しの元暗殺者。~無職のおっさんから始まるセカンドライフ~"
$ws.Cells.Item(25, 3).Value = 4
$ws.Cells.Item(26, 1).Value = 25
$ws.Cells.Item(26, 2).Value = "いつでも自宅に帰れる俺は、異世界で行商人をはじめました"
$ws.Cells.Item(26, 3).Value = 7
$ws.Cells.Item(27, 1).Value = 26
$ws.Cells.Item(27, 2).Value = "ROPPEN-六篇-"
$ws.Cells.Item(27, 3).Value = 10
$ws.Cells.Item(28, 1).Value = 27
$ws.Cells.Item(28, 2).Value = "ステータス・オール∞(インフィニティ) ∞使いの最強能力者、異世界を自由気ままに暮らします!"
$ws.Cells.Item(28, 3).Value = 11
$ws.Cells.Item(29, 1).Value = 28
$ws.Cells.Item(29, 2).Value = "彼女、お借りします"
$ws.Cells.Item(29, 3).Value = 43
$ws.Cells.Item(30, 1).Value = 29
$ws.Cells.Item(30, 2).Value = "春の嵐とモンスター"
$ws.Cells.Item(30, 3).Value = 9
$ws.Cells.Item(31, 1).Value = 30
$ws.Cells.Item(31, 2).Value = "たまのこしいれ ―アシガールEDO―"
$ws.Cells.Item(31, 3).Value = 4
$ws.Cells.Item(32, 1).Value = 31
$ws.Cells.Item(32, 2).Value = "落ちない汚れを僕は何と呼べばよかったのか"
$ws.Cells.Item(32, 3).Value = 1
$ws.Cells.Item(33, 1).Value = 32
$ws.Cells.Item(33, 2).Value = "おひとり様のナナイさん"
$ws.Cells.Item(33, 3).Value = 1
$ws.Cells.Item(34, 1).Value = 33
$ws.Cells.Item(34, 2).Value = "BLOOD THE LAST VAMPIRE 2000"
$ws.Cells.Item(34, 3).Value = 2
$ws.Cells.Item(35, 1).Value = 34
$ws.Cells.Item(35, 2).Value = "BLOOD THE LAST VAMPIRE 2000"
$ws.Cells.Item(35, 3).Value = 3
$ws.Cells.Item(36, 1).Value = 35
$ws.Cells.Item(36, 2).Value = "脳の髄まで愛してる 黄色い淫夢①"
$ws.Cells.Item(36, 3).Value = 1
$ws.Cells.Item(37, 1).Value = 36
$ws.Cells.Item(37, 2).Value = "復讐の同窓会"
$ws.Cells.Item(37, 3).Value = 1
$ws.Cells.Item(38, 1).Value = 37
$ws.Cells.Item(38, 2).Value = "小悪魔系天使は今日もいじわる"
$ws.Cells.Item(38, 3).Value = 1
$ws.Cells.Item(39, 1).Value = 38
$ws.Cells.Item(39, 2).Value = "1・2のアッホ!!"
$ws.Cells.Item(39, 3).Value = 2
$ws.Cells.Item(40, 1).Value = 39
$ws.Cells.Item(40, 2).Value = "転生したらスライムだった件~魔物の国の歩き方~"
$ws.Cells.Item(40, 3).Value = 1
$ws.Cells.Item(41, 1).Value = 40
$ws.Cells.Item(41, 2).Value = "超人X"
$ws.Cells.Item(41, 3).Value = 14
$ws.Cells.Item(42, 1).Value = 41
$ws.Cells.Item(42, 2).Value = "転生無敗の異世界賢者~ゲームのジョブで楽しいセカンドライフ~"
$ws.Cells.Item(42, 3).Value = 5
$ws.Cells.Item(43, 1).Value = 42
$ws.Cells.Item(43, 2).Value = "ふたりソロキャンプ"
$ws.Cells.Item(43, 3).Value = 22
$ws.Cells.Item(44, 1).Value = 43
$ws.Cells.Item(44, 2).Value = "デッドマウント・デスプレイ"
$ws.Cells.Item(44, 3).Value = 16
$ws.Cells.Item(45, 1).Value = 44
$ws.Cells.Item(45, 2).Value = "ドラフトキング"
$ws.Cells.Item(45, 3).Value = 23
$ws.Cells.Item(46, 1).Value = 45
$ws.Cells.Item(46, 2).Value = "ゆうべはお楽しみでしたね"
$ws.Cells.Item(46, 3).Value = 12
$ws.Cells.Item(47, 1).Value = 46
$ws.Cells.Item(47, 2).Value = "憂国のモリアーティ"
$ws.Cells.Item(47, 3).Value = 21
$ws.Cells.Item(48, 1).Value = 47
$ws.Cells.Item(48, 2).Value = "落ちない汚れを僕は何と呼べばよかったのか"
$ws.Cells.Item(48, 3).Value = 2
$ws.Cells.Item(49, 1).Value = 48
$ws.Cells.Item(49, 2).Value = "落ちない汚れを僕は何と呼べばよかったのか"
$ws.Cells.Item(49, 3).Value = 3
$ws.Cells.Item(50, 1).Value = 49
$ws.Cells.Item(50, 2).Value = "おひとり様のナナイさん"
$ws.Cells.Item(50, 3).Value = 2
$ws.Cells.Item(51, 1).Value = 50
$ws.Cells.Item(51, 2).Value = "おひとり様のナナイさん"
$ws.Cells.Item(51, 3).Value = 3
$ws.Cells.Item(52, 1).Value = 51
$ws.Cells.Item(52, 2).Value = "魔王の愛妻は愛されない"
$ws.Cells.Item(52, 3).Value = 1
$ws.Cells.Item(53, 1).Value = 52
$ws.Cells.Item(53, 2).Value = "山奥育ちの俺のゆるり異世界生活~もふもふと最強たちに可愛がられて、二度目の人生満喫中~"
$ws.Cells.Item(53, 3).Value = 1
$ws.Cells.Item(54, 1).Value = 53
$ws.Cells.Item(54, 2).Value = "1・2のアッホ!!"
$ws.Cells.Item(54, 3).Value = 3
$ws.Cells.Item(55, 1).Value = 54
$ws.Cells.Item(55, 2).Value = "もふもふと行く、腹ペコ料理人の絶品グルメライフ 第2話"
$ws.Cells.Item(55, 3).Value = 2
$ws.Cells.Item(56, 1).Value = 55
$ws.Cells.Item(56, 2).Value = "この恋、おくちにあいますか? ~優等生の白姫さんは問題児の俺と毎日キスしてる~"
$ws.Cells.Item(56, 3).Value = 1
$ws.Cells.Item(57, 1).Value = 56
$ws.Cells.Item(57, 2).Value = "この恋、おくちにあいますか? ~優等生の白姫さんは問題児の俺と毎日キスしてる~"
$ws.Cells.Item(57, 3).Value = 2
$ws.Cells.Item(58, 1).Value = 57
$ws.Cells.Item(58, 2).Value = "この恋、おくちにあいますか? ~優等生の白姫さんは問題児の俺と毎日キスしてる~"
$ws.Cells.Item(58, 3).Value = 3
$ws.Cells.Item(59, 1).Value = 58
$ws.Cells.Item(59, 2).Value = "嗤うサレ妻 復讐の好機は逃さない"
$ws.Cells.Item(59, 3).Value = 4
$ws.Cells.Item(60, 1).Value = 59
$ws.Cells.Item(60, 2).Value = "無能と追放された最弱魔法剣士、呪いが解けたので最強へ成り上がる6"
$ws.Cells.Item(60, 3).Value = 6
$ws.Cells.Item(61, 1).Value = 60
$ws.Cells.Item(61, 2).Value = "突撃!自衛官妻"
$ws.Cells.Item(61, 3).Value = 1
$ws.Cells.Item(62, 1).Value = 61
$ws.Cells.Item(62, 2).Value = "F REGENERATION 瑠璃"
$ws.Cells.Item(62, 3).Value = 1
$ws.Cells.Item(63, 1).Value = 62
$ws.Cells.Item(63, 2).Value = "風の騎士団"
$ws.Cells.Item(63, 3).Value = 1
$ws.Cells.Item(64, 1).Value = 63
$ws.Cells.Item(64, 2).Value = "Re-Tune ~あなたの人生チューニングします~"
$ws.Cells.Item(64, 3).Value = 1
$ws.Cells.Item(65, 1).Value = 64
$ws.Cells.Item(65, 2).Value = "天空の扉"
$ws.Cells.Item(65, 3).Value = 23
$ws.Cells.Item(66, 1).Value = 65
$ws.Cells.Item(66, 2).Value = "胚培養士(はいばいようし)ミズイロ~不妊治療のスペシャリスト~"
$ws.Cells.Item(66, 3).Value = 9
$ws.Cells.Item(67, 1).Value = 66
$ws.Cells.Item(67, 2).Value = "食料生成スキルを手に入れたので、異世界で商会を立ち上げようと思います:"
$ws.Cells.Item(67, 3).Value = 5
$ws.Cells.Item(68, 1).Value = 67
$ws.Cells.Item(68, 2).Value = "最後にひとつだけお願いしてもよろしいでしょうか11"
$ws.Cells.Item(68, 3).Value = 13
$ws.Cells.Item(69, 1).Value = 68
$ws.Cells.Item(69, 2).Value = "売れ残りの奴隷エルフを拾ったので、娘にすることにした"
$ws.Cells.Item(69, 3).Value = 3
$ws.Cells.Item(70, 1).Value = 69
$ws.Cells.Item(70, 2).Value = "シーカーズ~迷宮最強のおじさん、神配信者となる~"
$ws.Cells.Item(70, 3).Value = 5
$ws.Cells.Item(71, 1).Value = 70
$ws.Cells.Item(71, 2).Value = "異世界でスローライフを(願望)"
$ws.Cells.Item(71, 3).Value = 10
$ws.Cells.Item(72, 1).Value = 71
$ws.Cells.Item(72, 2).Value = "現実主義勇者の王国再建記XIV"
$ws.Cells.Item(72, 3).Value = 14
$ws.Cells.Item(73, 1).Value = 72
$ws.Cells.Item(73, 2).Value = "アフターゴッド"
$ws.Cells.Item(73, 3).Value = 10
$ws.Cells.Item(74, 1).Value = 73
$ws.Cells.Item(74, 2).Value = "黄昏流星群"
$ws.Cells.Item(74, 3).Value = 76
$ws.Cells.Item(75, 1).Value = 74
$ws.Cells.Item(75, 2).Value = "九条の大罪"
$ws.Cells.Item(75, 3).Value = 15
$ws.Cells.Item(76, 1).Value = 75
$ws.Cells.Item(76, 2).Value = "本好きの下剋上~司書になるためには手段を選んでいられません~第四部「貴族院の図書館を救いたい!11」"
$ws.Cells.Item(76, 3).Value = 11
$ws.Cells.Item(77, 1).Value = 76
$ws.Cells.Item(77, 2).Value = "カッコウの許嫁"
$ws.Cells.Item(77, 3).Value = 30
$ws.Cells.Item(78, 1).Value = 77
$ws.Cells.Item(78, 2).Value = "ガチアクタ"
$ws.Cells.Item(78, 3).Value = 17
$ws.Cells.Item(79, 1).Value = 78
$ws.Cells.Item(79, 2).Value = "虚構推理"
$ws.Cells.Item(79, 3).Value = 24
$ws.Cells.Item(80, 1).Value = 79
$ws.Cells.Item(80, 2).Value = "ワンパンマン"
$ws.Cells.Item(80, 3).Value = 34
$ws.Cells.Item(81, 1).Value = 80
$ws.Cells.Item(81, 2).Value = "ピンクとハバネロ"
$ws.Cells.Item(81, 3).Value = 14
$ws.Cells.Item(82, 1).Value = 81
$ws.Cells.Item(82, 2).Value = "旦那様の溺愛には黒い秘密がある"
$ws.Cells.Item(82, 3).Value = 1
$ws.Cells.Item(83, 1).Value = 82
$ws.Cells.Item(83, 2).Value = "雛名はシてあげたい!~アナタの復讐、引き受けます~"
$ws.Cells.Item(83, 3).Value = 1
$ws.Cells.Item(84, 1).Value = 83
$ws.Cells.Item(84, 2).Value = "目玉焼きにはソースか?醤油か?(フルカラー)"
$ws.Cells.Item(84, 3).Value = 1
$ws.Cells.Item(85, 1).Value = 84
$ws.Cells.Item(85, 2).Value = "復讐の刻~愛する息子が殺された~"
$ws.Cells.Item(85, 3).Value = 1
$ws.Cells.Item(86, 1).Value = 85
$ws.Cells.Item(86, 2).Value = "社長令嬢復讐日記"
$ws.Cells.Item(86, 3).Value = 1
$ws.Cells.Item(87, 1).Value = 86
$ws.Cells.Item(87, 2).Value = "無能力の私が次期当主のツガイになりました"
$ws.Cells.Item(87, 3).Value = 1
$ws.Cells.Item(88, 1).Value = 87
$ws.Cells.Item(88, 2).Value = "代わりの花嫁 ~愛する人と、姉の代わりに結婚します~"
$ws.Cells.Item(88, 3).Value = 1
$ws.Cells.Item(89, 1).Value = 88
$ws.Cells.Item(89, 2).Value = "うちの夫、やばくないですか?"
$ws.Cells.Item(89, 3).Value = 1
$ws.Cells.Item(90, 1).Value = 89
$ws.Cells.Item(90, 2).Value = "初恋は先生でした ~私を見つけてくれた人~"
$ws.Cells.Item(90, 3).Value = 1
$ws.Cells.Item(91, 1).Value = 90
$ws.Cells.Item(91, 2).Value = "捨てられ貴族の無人島のびのび開拓記~ようやく自由を手に入れたので、もふもふたちと気まぐれスローライフを満喫します~"
$ws.Cells.Item(91, 3).Value = 1
$ws.Cells.Item(92, 1).Value = 91
$ws.Cells.Item(92, 2).Value = "イジワル同居人は御曹司!?"
$ws.Cells.Item(92, 3).Value = 1
$ws.Cells.Item(93, 1).Value = 92
$ws.Cells.Item(93, 2).Value = "お見合い婚にも初夜は必要ですか?"
$ws.Cells.Item(93, 3).Value = 1
$ws.Cells.Item(94, 1).Value = 93
$ws.Cells.Item(94, 2).Value = "すばらしき新世界(フルカラー)"
$ws.Cells.Item(94, 3).Value = 9
$ws.Cells.Item(95, 1).Value = 94
$ws.Cells.Item(95, 2).Value = "すばらしき新世界(フルカラー)"
$ws.Cells.Item(95, 3).Value = 10
$ws.Cells.Item(96, 1).Value = 95
$ws.Cells.Item(96, 2).Value = "ホームセンターごと呼び出された私の大迷宮リノベーション!"
$ws.Cells.Item(96, 3).Value = 1
$ws.Cells.Item(97, 1).Value = 96
$ws.Cells.Item(97, 2).Value = "ホームセンターごと呼び出された私の大迷宮リノベーション!"
$ws.Cells.Item(97, 3).Value = 2
$ws.Cells.Item(98, 1).Value = 97
$ws.Cells.Item(98, 2).Value = "ホームセンターごと呼び出された私の大迷宮リノベーション!"
$ws.Cells.Item(98, 3).Value = 3
$ws.Cells.Item(99, 1).Value = 98
$ws.Cells.Item(99, 2).Value = "ダチョウ獣人のはちゃめちゃ無双"
$ws.Cells.Item(99, 3).Value = 1
$ws.Cells.Item(100, 1).Value = 99
$ws.Cells.Item(100, 2).Value = "ダチョウ獣人のはちゃめちゃ無双"
$ws.Cells.Item(100, 3).Value = 2
$ws.Cells.Item(101, 1).Value = 100
$ws.Cells.Item(101, 2).Value = "ダチョウ獣人のはちゃめちゃ無双"
$ws.Cells.Item(101, 3).Value = 3

# Apply highlighted-volume style (style index 2: light-yellow fill) to matching C cells via format copy,
# reusing existing style so no new style entries are created in styles.xml
$srcSheet.Range("C12").Copy()
$ws.Range("C8").PasteSpecial(-4122)
$ws.Range("C11:C12").PasteSpecial(-4122)
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C20:C22").PasteSpecial(-4122)
$ws.Range("C32:C40").PasteSpecial(-4122)
$ws.Range("C48:C58").PasteSpecial(-4122)
$ws.Range("C61:C64").PasteSpecial(-4122)
$ws.Range("C69").PasteSpecial(-4122)
$ws.Range("C82:C93").PasteSpecial(-4122)
$ws.Range("C96:C101").PasteSpecial(-4122)
